# Insert a new record at row 196 of the "Ajo" (garlic) price sheet.
# This pushes the existing rows 196-239 down to 197-240 (Excel keeps their
# data/styles intact automatically), and we then fill the freshly inserted
# row 196 with the new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 196, shifting everything below down.
$ws.Rows(196).Insert()

# Populate the new row 196 with the new weekly price record.
$ws.Range("A196").Value = 7
$ws.Range("B196").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C196").Value = 'Ñuble'
$ws.Range("D196").Value = 44754
$ws.Range("E196").Value = 16
$ws.Range("F196").Value = 100112003
$ws.Range("G196").Value = 'Ajo'
$ws.Range("H196").Value = 'Chino'
$ws.Range("I196").Value = 'Primera'
$ws.Range("J196").Value = 60
$ws.Range("K196").Value = 20000
$ws.Range("L196").Value = 21000
$ws.Range("M196").Value = 20500
$ws.Range("N196").Value = '$/caja 10 kilos'
$ws.Range("O196").Value = 'China'
$ws.Range("P196").Value = 2050
$ws.Range("Q196").Value = 10
$ws.Range("R196").Value = 'Hortaliza'
